$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.173.90"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.778.10"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.550"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.71"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.281"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0657"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D12").Value = "2.029.69"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.69%  "
$ws.Range("D14").Value = "1.776.65"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "34.138.42"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.79"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "0.0₃0738"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0511"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "1.441.13"
$ws.Range("E35").Value = "  -6.75%  "
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.626"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0187"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.890"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.05"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0508"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.80"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.931.99"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.998"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.65"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.52%  "
